$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.077.90'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '2.824.26'
$ws.Range("E3").Value = '  +8.86%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.42'
$ws.Range("E5").Value = '  +3.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '597.15'
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +4.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.193'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").Value = '2.824.96'
$ws.Range("E10").Value = '  +8.91%  '
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.370'
$ws.Range("E12").Value = '  +3.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.83'
$ws.Range("E13").Value = '  +2.83%  '
$ws.Range("D14").Value = '3.344.91'
$ws.Range("E14").Value = '  +7.47%  '
$ws.Range("D15").Value = '74.990.66'
$ws.Range("E15").Value = '  +2.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000188'
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.09'
$ws.Range("E17").Value = '  +5.19%  '
$ws.Range("D18").Value = '2.820.99'
$ws.Range("E18").Value = '  +8.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.09'
$ws.Range("E19").Value = '  +4.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  +5.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.75'
$ws.Range("E21").Value = '  +1.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.25'
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.79'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = '2.971.77'
$ws.Range("E27").Value = '  +8.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.18'
$ws.Range("E28").Value = '  +2.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  +4.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000103'
$ws.Range("E30").Value = '  +12.22%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '514.65'
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.39'
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.81'
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("E35").Value = '  +4.64%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.16'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.14'
$ws.Range("E38").Value = '  +5.95%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '183.54'
$ws.Range("E41").Value = '  +17.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.07'
$ws.Range("E43").Value = '  +5.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.341'
$ws.Range("E44").Value = '  +6.19%  '
$ws.Range("E45").Value = '  +2.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +5.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.99'
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0863'
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.34'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.569'
$ws.Range("E50").Value = '  +9.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.74'
$ws.Range("E51").Value = '  +4.65%  '
